$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 90: only the "Fecha" (date, column D) value changes
$ws.Range("D90").Value = 44463

# Rows 91-113: the date/volume/price block (D, J, K, L, M, P) shifts down by one row
# (each row now holds the values that used to belong to the row above it)
$ws.Range("D91").Value = 44251
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 900
$ws.Range("L91").Value = 1000
$ws.Range("M91").Value = 950
$ws.Range("P91").Value = 158
$ws.Range("D92").Value = 44433
$ws.Range("J92").Value = 3400
$ws.Range("K92").Value = 900
$ws.Range("L92").Value = 1000
$ws.Range("M92").Value = 950
$ws.Range("P92").Value = 158
$ws.Range("D93").Value = 44221
$ws.Range("J93").Value = 2600
$ws.Range("K93").Value = 1000
$ws.Range("L93").Value = 1100
$ws.Range("M93").Value = 1050
$ws.Range("P93").Value = 175
$ws.Range("D94").Value = 44316
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 900
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = 950
$ws.Range("P94").Value = 158
$ws.Range("D95").Value = 44279
$ws.Range("J95").Value = 3000
$ws.Range("K95").Value = 800
$ws.Range("L95").Value = 1000
$ws.Range("M95").Value = 900
$ws.Range("P95").Value = 150
$ws.Range("D96").Value = 44363
$ws.Range("J96").Value = 3360
$ws.Range("K96").Value = 900
$ws.Range("L96").Value = 1000
$ws.Range("M96").Value = 950
$ws.Range("P96").Value = 158
$ws.Range("D97").Value = 44277
$ws.Range("J97").Value = 2400
$ws.Range("K97").Value = 900
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = 950
$ws.Range("P97").Value = 158
$ws.Range("D98").Value = 44291
$ws.Range("J98").Value = 2600
$ws.Range("K98").Value = 900
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 950
$ws.Range("P98").Value = 158
$ws.Range("D99").Value = 44438
$ws.Range("J99").Value = 3100
$ws.Range("K99").Value = 900
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = 950
$ws.Range("P99").Value = 158
$ws.Range("D100").Value = 44372
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 900
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = 950
$ws.Range("P100").Value = 158
$ws.Range("D101").Value = 44286
$ws.Range("J101").Value = 3200
$ws.Range("K101").Value = 800
$ws.Range("L101").Value = 1000
$ws.Range("M101").Value = 900
$ws.Range("P101").Value = 150
$ws.Range("D102").Value = 44209
$ws.Range("J102").Value = 2700
$ws.Range("K102").Value = 1000
$ws.Range("L102").Value = 1100
$ws.Range("M102").Value = 1050
$ws.Range("P102").Value = 175
$ws.Range("D103").Value = 44356
$ws.Range("J103").Value = 3360
$ws.Range("K103").Value = 900
$ws.Range("L103").Value = 1000
$ws.Range("M103").Value = 950
$ws.Range("P103").Value = 158
$ws.Range("D104").Value = 44160
$ws.Range("J104").Value = 2700
$ws.Range("K104").Value = 800
$ws.Range("L104").Value = 1000
$ws.Range("M104").Value = 900
$ws.Range("P104").Value = 150
$ws.Range("D105").Value = 44351
$ws.Range("J105").Value = 2960
$ws.Range("K105").Value = 900
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 950
$ws.Range("P105").Value = 158
$ws.Range("D106").Value = 44365
$ws.Range("J106").Value = 2900
$ws.Range("K106").Value = 900
$ws.Range("L106").Value = 1000
$ws.Range("M106").Value = 950
$ws.Range("P106").Value = 158
$ws.Range("D107").Value = 44162
$ws.Range("J107").Value = 2400
$ws.Range("K107").Value = 800
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 900
$ws.Range("P107").Value = 150
$ws.Range("D108").Value = 44410
$ws.Range("J108").Value = 3200
$ws.Range("K108").Value = 900
$ws.Range("L108").Value = 1000
$ws.Range("M108").Value = 950
$ws.Range("P108").Value = 158
$ws.Range("D109").Value = 44244
$ws.Range("J109").Value = 2800
$ws.Range("K109").Value = 900
$ws.Range("L109").Value = 1000
$ws.Range("M109").Value = 950
$ws.Range("P109").Value = 158
$ws.Range("D110").Value = 44176
$ws.Range("J110").Value = 2400
$ws.Range("K110").Value = 800
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 900
$ws.Range("P110").Value = 150
$ws.Range("D111").Value = 44239
$ws.Range("J111").Value = 2600
$ws.Range("K111").Value = 1000
$ws.Range("L111").Value = 1100
$ws.Range("M111").Value = 1050
$ws.Range("P111").Value = 175
$ws.Range("D112").Value = 44358
$ws.Range("J112").Value = 3000
$ws.Range("K112").Value = 900
$ws.Range("L112").Value = 1000
$ws.Range("M112").Value = 950
$ws.Range("P112").Value = 158
$ws.Range("D113").Value = 44211
$ws.Range("J113").Value = 2400
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 1100
$ws.Range("M113").Value = 1050
$ws.Range("P113").Value = 175

# Row 114 is new: it receives the values that used to belong to row 113,
# together with the constant columns shared by the whole table
$ws.Range("A114").Value = 8
$ws.Range("B114").Value = "Terminal La Palmera de La Serena"
$ws.Range("C114").Value = "Coquimbo"
$ws.Range("D114").Value = 44323
$ws.Range("E114").Value = 4
$ws.Range("F114").Value = 100112037
$ws.Range("G114").Value = "Cebollín"
$ws.Range("H114").Value = "Sin especificar"
$ws.Range("I114").Value = "Primera"
$ws.Range("J114").Value = 2880
$ws.Range("K114").Value = 900
$ws.Range("L114").Value = 1000
$ws.Range("M114").Value = 950
$ws.Range("N114").Value = "$/paquete 6 unidades"
$ws.Range("O114").Value = "Provincia del Elquí"
$ws.Range("P114").Value = 158
$ws.Range("Q114").Value = 6
$ws.Range("R114").Value = "Hortaliza"

# Copy the date number-format style from D113 onto the new D114 cell
$ws.Range("D113").Copy()
$ws.Range("D114").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D114").Value = 44323
